$wb = $excel.ActiveWorkbook

# --- Overview sheet: the zh-cn / de-de status columns for the cb4dc240 row
#     move from "Ready for handoff" to "Handed back: in sync with en-US"
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn detail sheet: handback has completed for the cb4dc240 file,
#     so Status flips, the handback datetime is refreshed, and the stale
#     "old handback file" error is cleared.
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("K3").Value = "2016-08-30 06:51:28"
$zhcn.Range("P3").Value = ""
$zhcn.Columns.Item(16).ColumnWidth = 13.7470528738839

# --- de-de detail sheet: same handback update, with its own datetime
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K3").Value = "2016-08-30 06:51:35"
$dede.Range("P3").Value = ""
$dede.Columns.Item(16).ColumnWidth = 13.7470528738839
